{"js": "// Insert a new bulleted list item after the last paragraph in the document\n// (\"Separation of data from behaviours using Scriptable Objects\"), matching\n// the diff: a new ListParagraph-styled paragraph (numId 1, ilvl 0) with the\n// text \"Physics2D Settings Unchecking Queries Start in Colliders to disable\n// detecting itself.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// insertParagraph after the last paragraph inherits that paragraph's\n// formatting (style + list numbering), matching the surrounding bullet list.\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Physics2D Settings Unchecking Queries Start in Colliders to disable detecting itself.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item after the last paragraph in the document\n# (\"Separation of data from behaviours using Scriptable Objects\"), matching\n# the diff: a new ListParagraph-styled paragraph (numId 1, ilvl 0) with the\n# text \"Physics2D Settings Unchecking Queries Start in Colliders to disable\n# detecting itself.\"\n$d = $word.ActiveDocument\n\n$lastParaIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastParaIndex)\n\n# InsertParagraphAfter on the last paragraph's Range creates a new paragraph\n# that inherits the source paragraph's formatting (style + list numbering),\n# matching the surrounding bullet list.\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($lastParaIndex + 1)\n$newRange = $newPara.Range\n$newRange.Collapse(0)\n$newRange.Text = \"Physics2D Settings Unchecking Queries Start in Colliders to disable detecting itself.\"\n"}
